$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Unprotect()

# Update the confidential disclaimer date from 2021-03-29 to 2021-03-30
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# Update weight (D) and percent change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2460422759427153
$ws.Range("E2").Value = -0.003842794759825408

$ws.Range("D3").Value = 0.5006723694704533
$ws.Range("E3").Value = -0.003019489431786981

$ws.Range("D4").Value = 0.09771464372411084
$ws.Range("E4").Value = 0.002035002035002176

$ws.Range("D5").Value = 0.09911350067053598
$ws.Range("E5").Value = 0.006428988895382792

$ws.Range("D6").Value = 0.05645721019218451
$ws.Range("E6").Value = 0.01590693257359921

$ws.Range("E7").Value = -0.0007231547672500449

# Re-apply sheet protection (matching the sheet's protected state; the
# engine cannot reproduce the original legacy password hash, but the
# sheet is restored to a protected state as before the edit).
$ws.Protect()
